# Weekly data refresh for Hortaliza - Macroferia Regional de Talca - Zapallo italiano.
# A new week's record is prepended (date 44516) while every existing record shifts down
# one row (row N <- old row N-1), and the previously-last row becomes a new row 231.
# Column order: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Categoria ID,
# G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo, L Precio maximo,
# M Precio promedio ponderado, N Unidad de comercializacion, O Origen, P Precio $/Kg,
# Q Kg o Unidades, R Clasificacion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 219
$ws.Range("A219").Value = 5
$ws.Range("B219").Value = "Macroferia Regional de Talca"
$ws.Range("C219").Value = "Maule"
$ws.Range("D219").Value = 44516
$ws.Range("E219").Value = 7
$ws.Range("F219").Value = 100112032
$ws.Range("G219").Value = "Zapallo italiano"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 400
$ws.Range("K219").Value = 6000
$ws.Range("L219").Value = 6000
$ws.Range("M219").Value = 6000
$ws.Range("N219").Value = "`$/caja 60 unidades"
$ws.Range("O219").Value = "Región del Maule"
$ws.Range("P219").Value = 100
$ws.Range("Q219").Value = 60
$ws.Range("R219").Value = "Hortaliza"

# Row 220
$ws.Range("A220").Value = 5
$ws.Range("B220").Value = "Macroferia Regional de Talca"
$ws.Range("C220").Value = "Maule"
$ws.Range("D220").Value = 44270
$ws.Range("E220").Value = 7
$ws.Range("F220").Value = 100112032
$ws.Range("G220").Value = "Zapallo italiano"
$ws.Range("H220").Value = "Sin especificar"
$ws.Range("I220").Value = "Primera"
$ws.Range("J220").Value = 400
$ws.Range("K220").Value = 6000
$ws.Range("L220").Value = 6000
$ws.Range("M220").Value = 6000
$ws.Range("N220").Value = "`$/caja 60 unidades"
$ws.Range("O220").Value = "Región del Maule"
$ws.Range("P220").Value = 100
$ws.Range("Q220").Value = 60
$ws.Range("R220").Value = "Hortaliza"

# Row 221
$ws.Range("A221").Value = 5
$ws.Range("B221").Value = "Macroferia Regional de Talca"
$ws.Range("C221").Value = "Maule"
$ws.Range("D221").Value = 44295
$ws.Range("E221").Value = 7
$ws.Range("F221").Value = 100112032
$ws.Range("G221").Value = "Zapallo italiano"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 400
$ws.Range("K221").Value = 7000
$ws.Range("L221").Value = 7000
$ws.Range("M221").Value = 7000
$ws.Range("N221").Value = "`$/caja 60 unidades"
$ws.Range("O221").Value = "Región del Maule"
$ws.Range("P221").Value = 117
$ws.Range("Q221").Value = 60
$ws.Range("R221").Value = "Hortaliza"

# Row 222
$ws.Range("A222").Value = 5
$ws.Range("B222").Value = "Macroferia Regional de Talca"
$ws.Range("C222").Value = "Maule"
$ws.Range("D222").Value = 44217
$ws.Range("E222").Value = 7
$ws.Range("F222").Value = 100112032
$ws.Range("G222").Value = "Zapallo italiano"
$ws.Range("H222").Value = "Sin especificar"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 400
$ws.Range("K222").Value = 8000
$ws.Range("L222").Value = 8000
$ws.Range("M222").Value = 8000
$ws.Range("N222").Value = "`$/caja 60 unidades"
$ws.Range("O222").Value = "Región del Maule"
$ws.Range("P222").Value = 133
$ws.Range("Q222").Value = 60
$ws.Range("R222").Value = "Hortaliza"

# Row 223
$ws.Range("A223").Value = 5
$ws.Range("B223").Value = "Macroferia Regional de Talca"
$ws.Range("C223").Value = "Maule"
$ws.Range("D223").Value = 44509
$ws.Range("E223").Value = 7
$ws.Range("F223").Value = 100112032
$ws.Range("G223").Value = "Zapallo italiano"
$ws.Range("H223").Value = "Sin especificar"
$ws.Range("I223").Value = "Primera"
$ws.Range("J223").Value = 500
$ws.Range("K223").Value = 7000
$ws.Range("L223").Value = 7000
$ws.Range("M223").Value = 7000
$ws.Range("N223").Value = "`$/caja 60 unidades"
$ws.Range("O223").Value = "Región del Maule"
$ws.Range("P223").Value = 117
$ws.Range("Q223").Value = 60
$ws.Range("R223").Value = "Hortaliza"

# Row 224
$ws.Range("A224").Value = 5
$ws.Range("B224").Value = "Macroferia Regional de Talca"
$ws.Range("C224").Value = "Maule"
$ws.Range("D224").Value = 44421
$ws.Range("E224").Value = 7
$ws.Range("F224").Value = 100112032
$ws.Range("G224").Value = "Zapallo italiano"
$ws.Range("H224").Value = "Sin especificar"
$ws.Range("I224").Value = "Primera"
$ws.Range("J224").Value = 300
$ws.Range("K224").Value = 10000
$ws.Range("L224").Value = 10000
$ws.Range("M224").Value = 10000
$ws.Range("N224").Value = "`$/caja 50 unidades"
$ws.Range("O224").Value = "Región de Arica y Parinacota"
$ws.Range("P224").Value = 200
$ws.Range("Q224").Value = 50
$ws.Range("R224").Value = "Hortaliza"

# Row 225
$ws.Range("A225").Value = 5
$ws.Range("B225").Value = "Macroferia Regional de Talca"
$ws.Range("C225").Value = "Maule"
$ws.Range("D225").Value = 44244
$ws.Range("E225").Value = 7
$ws.Range("F225").Value = 100112032
$ws.Range("G225").Value = "Zapallo italiano"
$ws.Range("H225").Value = "Sin especificar"
$ws.Range("I225").Value = "Primera"
$ws.Range("J225").Value = 400
$ws.Range("K225").Value = 6000
$ws.Range("L225").Value = 6000
$ws.Range("M225").Value = 6000
$ws.Range("N225").Value = "`$/caja 60 unidades"
$ws.Range("O225").Value = "Región del Maule"
$ws.Range("P225").Value = 100
$ws.Range("Q225").Value = 60
$ws.Range("R225").Value = "Hortaliza"

# Row 226
$ws.Range("A226").Value = 5
$ws.Range("B226").Value = "Macroferia Regional de Talca"
$ws.Range("C226").Value = "Maule"
$ws.Range("D226").Value = 44307
$ws.Range("E226").Value = 7
$ws.Range("F226").Value = 100112032
$ws.Range("G226").Value = "Zapallo italiano"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 400
$ws.Range("K226").Value = 7000
$ws.Range("L226").Value = 7000
$ws.Range("M226").Value = 7000
$ws.Range("N226").Value = "`$/caja 60 unidades"
$ws.Range("O226").Value = "Región del Maule"
$ws.Range("P226").Value = 117
$ws.Range("Q226").Value = 60
$ws.Range("R226").Value = "Hortaliza"

# Row 227
$ws.Range("A227").Value = 5
$ws.Range("B227").Value = "Macroferia Regional de Talca"
$ws.Range("C227").Value = "Maule"
$ws.Range("D227").Value = 44273
$ws.Range("E227").Value = 7
$ws.Range("F227").Value = 100112032
$ws.Range("G227").Value = "Zapallo italiano"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 300
$ws.Range("K227").Value = 5000
$ws.Range("L227").Value = 5000
$ws.Range("M227").Value = 5000
$ws.Range("N227").Value = "`$/caja 60 unidades"
$ws.Range("O227").Value = "Región del Maule"
$ws.Range("P227").Value = 83
$ws.Range("Q227").Value = 60
$ws.Range("R227").Value = "Hortaliza"

# Row 228
$ws.Range("A228").Value = 5
$ws.Range("B228").Value = "Macroferia Regional de Talca"
$ws.Range("C228").Value = "Maule"
$ws.Range("D228").Value = 44433
$ws.Range("E228").Value = 7
$ws.Range("F228").Value = 100112032
$ws.Range("G228").Value = "Zapallo italiano"
$ws.Range("H228").Value = "Sin especificar"
$ws.Range("I228").Value = "Primera"
$ws.Range("J228").Value = 300
$ws.Range("K228").Value = 12000
$ws.Range("L228").Value = 12000
$ws.Range("M228").Value = 12000
$ws.Range("N228").Value = "`$/caja 50 unidades"
$ws.Range("O228").Value = "Región de Arica y Parinacota"
$ws.Range("P228").Value = 240
$ws.Range("Q228").Value = 50
$ws.Range("R228").Value = "Hortaliza"

# Row 229
$ws.Range("A229").Value = 5
$ws.Range("B229").Value = "Macroferia Regional de Talca"
$ws.Range("C229").Value = "Maule"
$ws.Range("D229").Value = 44302
$ws.Range("E229").Value = 7
$ws.Range("F229").Value = 100112032
$ws.Range("G229").Value = "Zapallo italiano"
$ws.Range("H229").Value = "Sin especificar"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 300
$ws.Range("K229").Value = 8000
$ws.Range("L229").Value = 8000
$ws.Range("M229").Value = 8000
$ws.Range("N229").Value = "`$/caja 60 unidades"
$ws.Range("O229").Value = "Región del Maule"
$ws.Range("P229").Value = 133
$ws.Range("Q229").Value = 60
$ws.Range("R229").Value = "Hortaliza"

# Row 230
$ws.Range("A230").Value = 5
$ws.Range("B230").Value = "Macroferia Regional de Talca"
$ws.Range("C230").Value = "Maule"
$ws.Range("D230").Value = 44179
$ws.Range("E230").Value = 7
$ws.Range("F230").Value = 100112032
$ws.Range("G230").Value = "Zapallo italiano"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 300
$ws.Range("K230").Value = 7000
$ws.Range("L230").Value = 7000
$ws.Range("M230").Value = 7000
$ws.Range("N230").Value = "`$/caja 60 unidades"
$ws.Range("O230").Value = "Región del Maule"
$ws.Range("P230").Value = 117
$ws.Range("Q230").Value = 60
$ws.Range("R230").Value = "Hortaliza"

# Row 231
$ws.Range("A231").Value = 5
$ws.Range("B231").Value = "Macroferia Regional de Talca"
$ws.Range("C231").Value = "Maule"
$ws.Range("D231").Value = 44491
$ws.Range("D231").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E231").Value = 7
$ws.Range("F231").Value = 100112032
$ws.Range("G231").Value = "Zapallo italiano"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 200
$ws.Range("K231").Value = 13000
$ws.Range("L231").Value = 13000
$ws.Range("M231").Value = 13000
$ws.Range("N231").Value = "`$/caja 60 unidades"
$ws.Range("O231").Value = "Región del Maule"
$ws.Range("P231").Value = 217
$ws.Range("Q231").Value = 60
$ws.Range("R231").Value = "Hortaliza"
